$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "match" table -> "round" table (id/label rename + column + type updates)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "round"
$ws.Range("E3").Value = "winner"            # table1 column "result" -> "winner"
$ws.Range("C4").Value = "INT"               # mid type INTEGER -> INT
$ws.Range("E4").Value = "ENUM"              # winner type roles(team) -> ENUM

# ---------------------------------------------------------------------------
# 2. "player" table - type column TEXT -> VARCHAR(30)
# ---------------------------------------------------------------------------
$ws.Range("I4").Value = "VARCHAR(30)"

# ---------------------------------------------------------------------------
# 3. "participates" table - foreign key annotation match(mid) -> round(mid)
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "round(mid)"

# ---------------------------------------------------------------------------
# 4. "roles" table - team type TEXT -> ENUM
# ---------------------------------------------------------------------------
$ws.Range("I9").Value = "ENUM"

# ---------------------------------------------------------------------------
# 5. "kills" table - match(mid) -> round(mid), time type TEXT -> TIME
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = "round(mid)"
$ws.Range("H14").Value = "TIME"

# ---------------------------------------------------------------------------
# 6. "damages" table - match(mid) -> round(mid), time TEXT -> TIME,
#    damage INTEGER -> INT
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = "round(mid)"
$ws.Range("H19").Value = "TIME"
$ws.Range("I19").Value = "INT"

# ---------------------------------------------------------------------------
# 7. New "configs" block + table (mirrors the 2-column "player" block)
# ---------------------------------------------------------------------------
$ws.Range("H2:I5").Copy()
$ws.Range("B22").PasteSpecial(-4122)   # xlPasteFormats, brings over the section styling

$ws.Range("B22").Value = "configs"
$ws.Range("B23").Value = " "
$ws.Range("C23").Value = "filename"
$ws.Range("B24").Value = "Type"
$ws.Range("C24").Value = "VARCHAR(30)"
$ws.Range("B25").Value = "Example"
$ws.Range("C25").Value = "config.log"

$ws.Range("B22:C22").Merge()

$cfgTable = $ws.ListObjects.Add(1, $ws.Range("B23:C25"), $null, 1)
$cfgTable.Name = "Tabelle28"

# Row 21 is the blank thick-bottom divider row preceding the new block,
# matching the spacer rows used elsewhere in the sheet (1, 6, 11, 16).
$ws.Rows("21").RowHeight = 15.75

# ---------------------------------------------------------------------------
# 8. Column widths widened to fit the new longer type strings
# ---------------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 12.166666666666666   # -> stored width 13
$ws.Columns("I:I").ColumnWidth = 12.166666666666666   # -> stored width 13

# ---------------------------------------------------------------------------
# 9. Selection matches the first merged header cell
# ---------------------------------------------------------------------------
$ws.Range("B2:F2").Select()
